$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = 2035
$ws.Range("B5").Value = 0.8099999999999999
$ws.Range("C5").Value = 0.15
$ws.Range("D5").Value = 0.04
